$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "45.426.54"
$ws.Range("E2").Value = "  +6.55%  "

$ws.Range("D3").Value = "2.388.18"
$ws.Range("E3").Value = "  +3.82%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "114.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "319.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.634"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.63%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("E9").Value = "  +3.69%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0930"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.09%  "

$ws.Range("E12").Value = "  +4.88%  "

$ws.Range("E13").Value = "  +2.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.53%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.73%  "

$ws.Range("D16").Value = "2.749.65"
$ws.Range("E16").Value = "  -1.10%  "

$ws.Range("D17").Value = "2.385.32"
$ws.Range("E17").Value = "  +4.33%  "

$ws.Range("D18").Value = "45.390.19"
$ws.Range("E18").Value = "  +6.61%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.70%  "

$ws.Range("E20").Value = "  +3.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.18%  "

$ws.Range("E23").Value = "  +3.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.00%  "

$ws.Range("E25").Value = "  +6.51%  "

$ws.Range("E26").Value = "  -0.65%  "

$ws.Range("E27").Value = "  +5.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "39.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0981"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +14.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "172.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.44%  "

$ws.Range("E34").Value = "  +10.77%  "

$ws.Range("E35").Value = "  +2.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.16%  "

$ws.Range("E37").Value = "  +6.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.86%  "

$ws.Range("E40").Value = "  +4.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.242"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.93%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "71.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.80%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +14.50%  "

$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +13.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "116.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.75%  "

$ws.Range("E51").Value = "  +10.66%  "
